$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5,D6,D10,D12,D14,D19,D20,D22,D23,D24,D25,D26,D29,D30,D31,D34,D35,D36,D37,D38,D40,D41,D43,D44,D47,D50,D51").NumberFormat = "@"

$ws.Range('D2').Value = '63.423.10'
$ws.Range('E2').Value = '  +4.56%  '
$ws.Range('D3').Value = '3.056.57'
$ws.Range('E3').Value = '  +2.77%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').Value = '549.03'
$ws.Range('E5').Value = '  +4.69%  '
$ws.Range('D6').Value = '139.34'
$ws.Range('E6').Value = '  +7.16%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '3.050.72'
$ws.Range('E8').Value = '  +2.71%  '
$ws.Range('E9').Value = '  +2.57%  '
$ws.Range('D10').Value = '6.27'
$ws.Range('E10').Value = '  +3.14%  '
$ws.Range('E11').Value = '  +1.25%  '
$ws.Range('D12').Value = '0.451'
$ws.Range('E12').Value = '  +3.40%  '
$ws.Range('E13').Value = '  +4.45%  '
$ws.Range('D14').Value = '34.65'
$ws.Range('E14').Value = '  +4.77%  '
$ws.Range('D15').Value = '3.560.99'
$ws.Range('E15').Value = '  +3.26%  '
$ws.Range('D16').Value = '63.487.83'
$ws.Range('E16').Value = '  +4.76%  '
$ws.Range('D17').Value = '3.059.29'
$ws.Range('E17').Value = '  +3.22%  '
$ws.Range('E18').Value = '  -0.82%  '
$ws.Range('D19').Value = '6.70'
$ws.Range('E19').Value = '  +3.80%  '
$ws.Range('D20').Value = '478.58'
$ws.Range('E20').Value = '  +5.19%  '
$ws.Range('E21').Value = '  +4.10%  '
$ws.Range('D22').Value = '0.674'
$ws.Range('E22').Value = '  +1.61%  '
$ws.Range('D23').Value = '7.16'
$ws.Range('E23').Value = '  +5.82%  '
$ws.Range('D24').Value = '80.96'
$ws.Range('E24').Value = '  +3.88%  '
$ws.Range('D25').Value = '12.45'
$ws.Range('E25').Value = '  +6.25%  '
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.13%  '
$ws.Range('E27').Value = '  +4.96%  '
$ws.Range('E28').Value = '  +4.50%  '
$ws.Range('B29').Value = 'FirstDigitalUSD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.32%  '
$ws.Range('B30').Value = 'ImmutableX'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D30').Value = '1.98'
$ws.Range('E30').Value = '  +8.27%  '
$ws.Range('D31').Value = '25.93'
$ws.Range('E31').Value = '  +3.66%  '
$ws.Range('E32').Value = '  +1.23%  '
$ws.Range('E33').Value = '  +7.68%  '
$ws.Range('D34').Value = '5.64'
$ws.Range('E34').Value = '  +7.24%  '
$ws.Range('D35').Value = '55.92'
$ws.Range('E35').Value = '  +2.01%  '
$ws.Range('D36').Value = '5.96'
$ws.Range('E36').Value = '  +4.25%  '
$ws.Range('D37').Value = '464.10'
$ws.Range('E37').Value = '  +3.07%  '
$ws.Range('D38').Value = '0.0812'
$ws.Range('E38').Value = '  +5.38%  '
$ws.Range('D39').Value = '3.123.15'
$ws.Range('E39').Value = '  -1.33%  '
$ws.Range('D40').Value = '0.0393'
$ws.Range('E40').Value = '  +4.46%  '
$ws.Range('D41').Value = '0.118'
$ws.Range('E41').Value = '  +3.19%  '
$ws.Range('E42').Value = '  +2.80%  '
$ws.Range('D43').Value = '2.57'
$ws.Range('E43').Value = '  +7.54%  '
$ws.Range('D44').Value = '27.95'
$ws.Range('E44').Value = '  +11.41%  '
$ws.Range('E45').Value = '  +3.97%  '
$ws.Range('E46').Value = '  -0.12%  '
$ws.Range('D47').Value = '2.04'
$ws.Range('E47').Value = '  +5.87%  '
$ws.Range('E48').Value = '  +1.54%  '
$ws.Range('D49').Value = '0.0₃0508'
$ws.Range('E49').Value = '  +3.36%  '
$ws.Range('D50').Value = '115.69'
$ws.Range('E50').Value = '  -1.54%  '
$ws.Range('D51').Value = '2.05'
$ws.Range('E51').Value = '  +6.16%  '
